# Commit: "update file with jgit"
# Change the E8 cell on the "Rules" sheet from "Good Morning" to "GIT UPDATE",
# and leave the selection on that cell (matches the <selection activeCell="E8".../>
# recorded in the saved sheetView).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
